$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C3").Value = 2650.7782999999999
$ws.Range("C4").Value = 1894.2938999999999
$ws.Range("C5").Value = 16164.9395
$ws.Range("C7").Value = 120.9529

$ws.Range("C13").Select()
